$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BB1: header date, copy style (border+date format) from BA1, then set value ---
$ws.Range("BA1").Copy($ws.Range("BB1"))
$ws.Range("BB1").Value = 45986

# --- BB2:BB82: new forecast column, plain numeric (no special style) ---
$ws.Range("BB2").Value = 2.672233110627005
$ws.Range("BB3").Value = -0.2486584437591262
$ws.Range("BB4").Value = 2.48444986756347
$ws.Range("BB5").Value = 1.924160321525846
$ws.Range("BB6").Value = 1.36313906046999
$ws.Range("BB7").Value = -4.856442119446953
$ws.Range("BB8").Value = 1.575066711296429
$ws.Range("BB9").Value = 1.814138146527952
$ws.Range("BB10").Value = 0.4238544839342779
$ws.Range("BB11").Value = -0.1584284749490763
$ws.Range("BB12").Value = 0.4882046181590169
$ws.Range("BB13").Value = 1.53878081519332
$ws.Range("BB14").Value = 0.8692906535860487
$ws.Range("BB15").Value = 3.11595497587993
$ws.Range("BB16").Value = 0.3906876693375665
$ws.Range("BB17").Value = 1.221836833579857
$ws.Range("BB18").Value = 1.653194230429179
$ws.Range("BB19").Value = -1.198649426118308
$ws.Range("BB20").Value = 0.7741673943688596
$ws.Range("BB21").Value = 0.8244429407371285
$ws.Range("BB22").Value = -0.3344927361763723
$ws.Range("BB23").Value = 0.7295104243506501
$ws.Range("BB24").Value = 0.03704431742310987
$ws.Range("BB25").Value = 0.2220551683158618
$ws.Range("BB26").Value = 0.1661617031019347
$ws.Range("BB27").Value = 1.924586088852507
$ws.Range("BB28").Value = 0.5703441042800677
$ws.Range("BB29").Value = 0.5671096279522487
$ws.Range("BB30").Value = 0.6731221874372437
$ws.Range("BB31").Value = 0.8295169162459786
$ws.Range("BB32").Value = 0.363674885967896
$ws.Range("BB33").Value = 0.6363538952886927
$ws.Range("BB34").Value = 0.4303070273019074
$ws.Range("BB35").Value = 0.907245662456674
$ws.Range("BB36").Value = 0.7278878628511336
$ws.Range("BB37").Value = 0.5936080878907575
$ws.Range("BB38").Value = 0.3249989166702818
$ws.Range("BB39").Value = 2.043550613228959
$ws.Range("BB40").Value = 0.867287375484608
$ws.Range("BB41").Value = 0.6331942894404392
$ws.Range("BB42").Value = -0.05439614307451279
$ws.Range("BB43").Value = 0.1739459843577862
$ws.Range("BB44").Value = 1.255538557350434
$ws.Range("BB45").Value = 1.090483027535811
$ws.Range("BB46").Value = 1.226659036647675
$ws.Range("BB47").Value = -0.5
$ws.Range("BB48").Value = 0.9
$ws.Range("BB49").Value = 1
$ws.Range("BB50").Value = 1.1
$ws.Range("BB51").Value = -2.477834671711193
$ws.Range("BB52").Value = 0.2428240279789122
$ws.Range("BB53").Value = 0.8650544612728055
$ws.Range("BB54").Value = 0.05143518179183104
$ws.Range("BB55").Value = -2.588552528306963
$ws.Range("BB56").Value = 1.067142397791443
$ws.Range("BB57").Value = 0.9403228036019016
$ws.Range("BB58").Value = 1.348411706012428
$ws.Range("BB59").Value = -0.1190215178375666
$ws.Range("BB60").Value = 0.7996487817115536
$ws.Range("BB61").Value = 1.039074166251879
$ws.Range("BB62").Value = 0.6958556561364588
$ws.Range("BB63").Value = -1.636815679601384
$ws.Range("BB64").Value = -0.04357278727286484
$ws.Range("BB65").Value = 0.3050453114869214
$ws.Range("BB66").Value = 0.7386451510207621
$ws.Range("BB67").Value = 1.021259612058628
$ws.Range("BB68").Value = 0.5487159577757694
$ws.Range("BB69").Value = 0.6800497182067176
$ws.Range("BB70").Value = 0.1667761162031525
$ws.Range("BB71").Value = 0.5709980498538272
$ws.Range("BB72").Value = 1.176666004305858
$ws.Range("BB73").Value = 0.8783323788356512
$ws.Range("BB74").Value = 0.6042915512474423
$ws.Range("BB75").Value = 0.6042915512474423
$ws.Range("BB76").Value = 0.6042915512474423
$ws.Range("BB77").Value = 0.6042915512474423
$ws.Range("BB78").Value = 0.6042915512474423
$ws.Range("BB79").Value = 0.6042915512474423
$ws.Range("BB80").Value = 0.6042915512474423
$ws.Range("BB81").Value = 0.6042915512474423
$ws.Range("BB82").Value = 0.6042915512474423

# --- Row 83: new row. A83 date (style like A82), BB83 new forecast value ---
$ws.Range("A82").Copy($ws.Range("A83"))
$ws.Range("A83").Value = 46934
$ws.Range("BB83").Value = 0.6042915512474423
